$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.38601601592701229
$ws.Range("A2").Value = -0.0099999992126775794
$ws.Range("A3").Value = -0.085300523092787017
$ws.Range("A4").Value = -0.011999999795982319
$ws.Range("A5").Value = -0.0059999992045165484
$ws.Range("A6").Value = -0.0059999991844712497
$ws.Range("A7").Value = -0.019999999045486661
$ws.Range("A8").Value = -0.019999999047095152
$ws.Range("A9").Value = -0.005999999192123795
$ws.Range("A10").Value = 0.057755277470690203
$ws.Range("A11").Value = 0.031076518955948984
$ws.Range("A12").Value = -0.0059999991965042909
$ws.Range("A13").Value = -0.0059999991846195755
$ws.Range("A14").Value = -0.011999999121620419
$ws.Range("A15").Value = -0.0059999991786670037
$ws.Range("A16").Value = -0.0059999991762760274
$ws.Range("A17").Value = -0.0059999991730972368
$ws.Range("A18").Value = -0.0089999991426070736
$ws.Range("A19").Value = -0.0089999992229365944
$ws.Range("A20").Value = -0.0089999992158134035
$ws.Range("A21").Value = -0.0089999992147538066
$ws.Range("A22").Value = -0.0089999992140503693
$ws.Range("A23").Value = -0.0089999991854199379
$ws.Range("A24").Value = -0.041999998844443276
$ws.Range("A25").Value = -0.041999998837981778
$ws.Range("A26").Value = -0.059479742224826992
$ws.Range("A27").Value = -0.005999999178105675
$ws.Range("A28").Value = -0.0059999991671721986
$ws.Range("A29").Value = -0.011999999100988035
$ws.Range("A30").Value = -0.019999999018347481
$ws.Range("A31").Value = -0.014999999065265612
$ws.Range("A32").Value = 0.0043317948749441726
$ws.Range("A33").Value = -0.0059999991544712472
